$d = $word.ActiveDocument

# Locate "Manuela Ojeda Ojeda.  Cód.:" and append the missing student ID
# (" 201814476") right after it, matching the sibling entries above it
# (e.g. "Isabella Mendez Pedraza.  Cód.: 201814239").
$r = $d.Content
$found = $r.Find.Execute("Manuela Ojeda Ojeda.  Cód.:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r.Collapse(0)   # wdCollapseEnd
    $r.InsertAfter(" 201814476")
} else {
    Write-Host "WARNING: target text not found"
}
